# Update cryptos list (Wed Nov 22 23:02:31 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '37.484.10'
$ws.Range("E2").Value = '  +2.82%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '2.065.96'
$ws.Range("E3").Value = '  +6.28%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.15%  '

# Row 5 - BNB
$cell = $ws.Cells.Item(5, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '236.44'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  +5.00%  '

# Row 6 - XRP
$ws.Range("E6").Value = '  +5.04%  '

# Row 7 - Solana
$cell = $ws.Cells.Item(7, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '58.21'
$cell.Style = $origStyle
$ws.Range("E7").Value = '  +11.73%  '

# Row 9 - Cardano
$ws.Range("E9").Value = '  +6.36%  '

# Row 10 - OKB
$cell = $ws.Cells.Item(10, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '57.76'
$cell.Style = $origStyle
$ws.Range("E10").Value = '  +2.05%  '

# Row 11 - Dogecoin
$ws.Range("E11").Value = '  +5.90%  '

# Row 12 - TRON
$ws.Range("E12").Value = '  +5.24%  '

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '2.373.47'
$ws.Range("E13").Value = '  +5.92%  '

# Row 14 - Chainlink
$cell = $ws.Cells.Item(14, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '14.30'
$cell.Style = $origStyle
$ws.Range("E14").Value = '  +6.12%  '

# Row 15 - Avalanche
$cell = $ws.Cells.Item(15, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '20.91'
$cell.Style = $origStyle
$ws.Range("E15").Value = '  +8.30%  '

# Row 16 - Polygon
$ws.Range("E16").Value = '  +6.61%  '

# Row 17 - Polkadot
$ws.Range("E17").Value = '  +6.27%  '

# Row 18 - WrappedEther
$ws.Range("D18").Value = '2.065.18'
$ws.Range("E18").Value = '  +5.84%  '

# Row 19 - WrappedBTC
$ws.Range("D19").Value = '37.626.08'
$ws.Range("E19").Value = '  +3.49%  '

# Row 20 - Uniswap
$ws.Range("E20").Value = '  +25.16%  '

# Row 21 - Litecoin
$cell = $ws.Cells.Item(21, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '68.52'
$cell.Style = $origStyle
$ws.Range("E21").Value = '  +3.86%  '

# Row 22 - ShibaInu
$ws.Range("D22").Value = '0.0₃0810'
$ws.Range("E22").Value = '  +4.54%  '

# Row 23 - BitcoinCash
$cell = $ws.Cells.Item(23, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '224.85'
$cell.Style = $origStyle
$ws.Range("E23").Value = '  +3.97%  '

# Row 24 - Dai
$cell = $ws.Cells.Item(24, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = $origStyle
$ws.Range("E24").Value = '  +0.14%  '

# Row 25 - PancakeSwap
$ws.Range("E25").Value = '  +8.73%  '

# Row 26 - Toncoin
$ws.Range("E26").Value = '  +3.89%  '

# Row 27 - Monero
$cell = $ws.Cells.Item(27, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '162.80'
$cell.Style = $origStyle
$ws.Range("E27").Value = '  +2.62%  '

# Row 28 - Cosmos
$cell = $ws.Cells.Item(28, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.85'
$cell.Style = $origStyle
$ws.Range("E28").Value = '  +6.48%  '

# Row 29 - Kaspa
$cell = $ws.Cells.Item(29, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.129'
$cell.Style = $origStyle
$ws.Range("E29").Value = '  +9.67%  '

# Row 30 - ImmutableX
$ws.Range("E30").Value = '  +10.45%  '

# Row 31 - EthereumClassic
$cell = $ws.Cells.Item(31, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '19.24'
$cell.Style = $origStyle
$ws.Range("E31").Value = '  +5.04%  '

# Row 32 - Stellar
$ws.Range("E32").Value = '  +4.24%  '

# Row 33 - LidoDAOToken
$ws.Range("E33").Value = '  +18.42%  '

# Row 34 - Filecoin
$ws.Range("E34").Value = '  +6.91%  '

# Row 35 - Hedera
$ws.Range("E35").Value = '  +6.90%  '

# Row 36 - InternetComputer(DFINITY)
$cell = $ws.Cells.Item(36, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.45'
$cell.Style = $origStyle
$ws.Range("E36").Value = '  +8.99%  '

# Row 37 - BinanceUSD->WEMIXToken
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell = $ws.Cells.Item(37, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.80'
$cell.Style = $origStyle
$ws.Range("E37").Value = '  +1.16%  '

# Row 38 - WEMIXToken->BinanceUSD
$ws.Range("B38").Value = 'BinanceUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Cells.Item(38, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = $origStyle
$ws.Range("E38").Value = '  -0.18%  '

# Row 39 - RenderToken
$ws.Range("E39").Value = '  +7.27%  '

# Row 40 - THORChain
$ws.Range("E40").Value = '  +17.69%  '

# Row 41 - HuobiToken
$ws.Range("E41").Value = '  -1.96%  '

# Row 42 - FTXToken
$ws.Range("E42").Value = '  +32.82%  '

# Row 43 - Cronos
$cell = $ws.Cells.Item(43, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0954'
$cell.Style = $origStyle
$ws.Range("E43").Value = '  +11.57%  '

# Row 44 - Maker
$ws.Range("D44").Value = '1.470.12'
$ws.Range("E44").Value = '  +6.14%  '

# Row 45 - Aave
$cell = $ws.Cells.Item(45, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '95.52'
$cell.Style = $origStyle
$ws.Range("E45").Value = '  +12.76%  '

# Row 46 - VeChain
$ws.Range("E46").Value = '  +7.91%  '

# Row 47 - InjectiveProtocol
$ws.Range("E47").Value = '  +12.42%  '

# Row 48 - TrustWalletToken
$ws.Range("E48").Value = '  +7.63%  '

# Row 49 - FraxShare->ARBITRUM
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Cells.Item(49, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.02'
$cell.Style = $origStyle
$ws.Range("E49").Value = '  +6.63%  '

# Row 50 - ARBITRUM->FraxShare
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Cells.Item(50, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.27'
$cell.Style = $origStyle
$ws.Range("E50").Value = '  +9.82%  '

# Row 51 - MXToken
$ws.Range("E51").Value = '  +3.03%  '
